# Add a new "2022-Q1" sheet (holdings detail) positioned right before the
# "总计" (total/summary) sheet, and insert a corresponding new first data
# row ("2022-Q1", 10, 33.27) into the "总计" sheet.

$wb = $excel.ActiveWorkbook
$totalSheetName = "总计"

# ---------------------------------------------------------------------------
# 1. Locate the existing "总计" sheet (always the last sheet) and insert a
#    brand-new worksheet immediately before it.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item($totalSheetName)
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# NOTE: adding a new sheet "before" $totalSheet can leave the $totalSheet
# handle pointing at the wrong sheet; re-fetch it fresh by name to be safe.
$totalSheet = $wb.Worksheets.Item($totalSheetName)

# A worksheet that already uses the shared "header/index" style (s="2") so we
# can copy its formatting onto the new sheet's header row / index column.
$styleSource = $wb.Worksheets.Item("2021-Q2").Range("A2")

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q1" sheet with the fund holdings detail.
# ---------------------------------------------------------------------------
$ws = $newSheet

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = [char]([int][char]'B' + $i)
    $ws.Range($col + "1").Value = $headers[$i]
}

$rows = @(
    @("008903", "广发科技先锋混合",          "149.48", "94.86", "6.18", "9.2379", 10),
    @("005911", "广发双擎升级混合A",          "139.02", "94.60", "6.17", "8.5775", 7),
    @("162703", "广发小盘成长混合(LOF)A",     "98.71",  "94.88", "7.85", "7.7487", 5),
    @("002939", "广发创新升级灵活配置混合",    "94.55",  "94.85", "5.95", "5.6257", 9),
    @("004854", "广发中证全指汽车指数A",       "22.01",  "94.43", "5.02", "1.1049", 6),
    @("009132", "广发小盘成长混合(LOF)C",     "5.31",   "94.88", "7.85", "0.4168", 5),
    @("004855", "广发中证全指汽车指数C",       "6.11",   "94.43", "5.02", "0.3067", 6),
    @("009314", "广发双擎升级混合C",          "3.63",   "94.60", "6.17", "0.2240", 7),
    @("006692", "金信消费升级股票A",          "0.62",   "94.05", "3.26", "0.0202", 10),
    @("006693", "金信消费升级股票C",          "0.20",   "94.05", "3.26", "0.0065", 10)
)

$lastRow = 1 + $rows.Length

# Text-like numeric-looking values (fund code, scale, positions, ratios, nav)
# need to be forced to text *before* assignment, otherwise Excel silently
# reinterprets them as numbers and e.g. drops leading zeros.
$textRange = $ws.Range("B2:G" + $lastRow)
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Range("B" + $r).Value = $row[0]
    $ws.Range("C" + $r).Value = $row[1]
    $ws.Range("D" + $r).Value = $row[2]
    $ws.Range("E" + $r).Value = $row[3]
    $ws.Range("F" + $r).Value = $row[4]
    $ws.Range("G" + $r).Value = $row[5]
    $ws.Range("H" + $r).Value = $row[6]
    $ws.Range("A" + $r).Value = $i
}

# Drop the temporary "@" number format again (content already committed as
# text) so the plain data cells end up with the default (no) style, exactly
# like the other quarter sheets.
$textRange.Style = "Normal"

# Apply the shared header/index style to the header row and the index column.
$styleSource.Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$ws.Range("A2:A" + $lastRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Insert a new first data row into the "总计" sheet for 2022-Q1, pushing
#    the previously existing rows down by one.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item($totalSheetName)
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 10
$totalSheet.Range("D2").Value = 33.27
$totalSheet.Range("A2").Value = 0

# Fix up the index column (A) below the newly inserted row: it should keep
# counting 0,1,2,... rather than the shifted original values.
$totalLastRow = $totalSheet.UsedRange.Rows.Count
for ($r = 3; $r -le $totalLastRow; $r++) {
    $totalSheet.Range("A" + $r).Value = $r - 2
}

# The inserted row copied formatting from the row above (the header), which
# is not what we want for a plain data row; restore plain formatting for the
# B2:D2 cells while keeping A2's index style.
$plainSource = $wb.Worksheets.Item("2021-Q2").Range("B2")
$plainSource.Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$indexStyleSource = $wb.Worksheets.Item("2021-Q2").Range("A2")
$indexStyleSource.Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
